$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the multiple runs in the "git log -until=..." paragraph into one.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "git log –until=”2024-19-09” = commits done till “date”",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "git log –until=”2024-19-09” = commits done till “date”", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the multiple runs in the "after making any changes to file in
#    github, to push the changes to local" paragraph into one.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "after making any changes to file in github, to push the changes to local",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "after making any changes to file in github, to push the changes to local", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert the new branching instructions right after the "git pull"
#    paragraph (and before the trailing blank paragraphs).
# ---------------------------------------------------------------------------
function Insert-ParaAfter {
    param(
        [int]$Index,
        [string]$Text
    )
    $r = $d.Paragraphs.Item($Index).Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    if ($Text) {
        $newRange = $d.Paragraphs.Item($Index + 1).Range
        $newRange.InsertAfter($Text)
    }
}

# Find the "git pull" paragraph that follows the just-merged paragraph.
$gitPullIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "git pull`r") {
        $gitPullIndex = $i
    }
}

$idx = $gitPullIndex
Insert-ParaAfter $idx ""                                                ; $idx++
Insert-ParaAfter $idx "To create a new branch from an existing branch " ; $idx++
Insert-ParaAfter $idx "git branch branch_name"                          ; $idx++
Insert-ParaAfter $idx ""                                                ; $idx++
Insert-ParaAfter $idx "to move from current branch to new branch"       ; $idx++
Insert-ParaAfter $idx "git checkout branch_name"                        ; $idx++
